$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.04383033333333
$ws.Range("H2").Value = 42.131491
$ws.Range("I2").Value = 0.158930310642385
$ws.Range("J2").Value = 0.158930310642385
$ws.Range("M2").Value = 0.01989833333333333
$ws.Range("N2").Value = 0.059695
$ws.Range("O2").Value = 0.5455534129646046
$ws.Range("P2").Value = 0.5455534129646046
$ws.Range("Q2").Value = 0.2794488172494444
$ws.Range("R2").Value = 2.515039355245
$ws.Range("S2").Value = 0.08670497339447795
$ws.Range("T2").Value = 0.08670497339447794
$ws.Range("G3").Value = 14.04383033333333
$ws.Range("H3").Value = 42.131491
$ws.Range("I3").Value = 0.158930310642385
$ws.Range("J3").Value = 0.158930310642385
$ws.Range("O3").Value = 0.1997148627777118
$ws.Range("P3").Value = 0.1997148627777118
$ws.Range("Q3").Value = 0.1022999414247778
$ws.Range("R3").Value = 0.920699472823
$ws.Range("S3").Value = 0.03174074518116302
$ws.Range("T3").Value = 0.03174074518116302
$ws.Range("G4").Value = 14.04383033333333
$ws.Range("H4").Value = 42.131491
$ws.Range("I4").Value = 0.158930310642385
$ws.Range("J4").Value = 0.158930310642385
$ws.Range("M4").Value = 0.009290999999999999
$ws.Range("N4").Value = 0.027873
$ws.Range("O4").Value = 0.2547317242576836
$ws.Range("P4").Value = 0.2547317242576836
$ws.Range("Q4").Value = 0.130481227627
$ws.Range("R4").Value = 1.174331048643
$ws.Range("S4").Value = 0.04048459206674401
$ws.Range("T4").Value = 0.04048459206674401
$ws.Range("I5").Value = 0.3128977537755355
$ws.Range("J5").Value = 0.3128977537755354
$ws.Range("M5").Value = 0.01989833333333333
$ws.Range("N5").Value = 0.059695
$ws.Range("O5").Value = 0.5455534129646046
$ws.Range("P5").Value = 0.5455534129646046
$ws.Range("Q5").Value = 0.5501713729694445
$ws.Range("R5").Value = 4.951542356725
$ws.Range("S5").Value = 0.1707024374812019
$ws.Range("T5").Value = 0.1707024374812018
$ws.Range("I6").Value = 0.3128977537755355
$ws.Range("J6").Value = 0.3128977537755354
$ws.Range("O6").Value = 0.1997148627777118
$ws.Range("P6").Value = 0.1997148627777118
$ws.Range("S6").Value = 0.06249033195873532
$ws.Range("T6").Value = 0.06249033195873531
$ws.Range("I7").Value = 0.3128977537755355
$ws.Range("J7").Value = 0.3128977537755354
$ws.Range("M7").Value = 0.009290999999999999
$ws.Range("N7").Value = 0.027873
$ws.Range("O7").Value = 0.2547317242576836
$ws.Range("P7").Value = 0.2547317242576836
$ws.Range("Q7").Value = 0.256887958435
$ws.Range("R7").Value = 2.311991625915
$ws.Range("S7").Value = 0.07970498433559828
$ws.Range("T7").Value = 0.07970498433559828
$ws.Range("G8").Value = 25.89747433333334
$ws.Range("H8").Value = 77.69242300000001
$ws.Range("I8").Value = 0.2930748622675039
$ws.Range("J8").Value = 0.2930748622675038
$ws.Range("M8").Value = 0.01989833333333333
$ws.Range("N8").Value = 0.059695
$ws.Range("O8").Value = 0.5455534129646046
$ws.Range("P8").Value = 0.5455534129646046
$ws.Range("Q8").Value = 0.5153165767761112
$ws.Range("R8").Value = 4.637849190985
$ws.Range("S8").Value = 0.1598879913641682
$ws.Range("T8").Value = 0.1598879913641681
$ws.Range("G9").Value = 25.89747433333334
$ws.Range("H9").Value = 77.69242300000001
$ws.Range("I9").Value = 0.2930748622675039
$ws.Range("J9").Value = 0.2930748622675038
$ws.Range("O9").Value = 0.1997148627777118
$ws.Range("P9").Value = 0.1997148627777118
$ws.Range("Q9").Value = 0.1886458355354445
$ws.Range("R9").Value = 1.697812519819
$ws.Range("S9").Value = 0.05853140590135132
$ws.Range("T9").Value = 0.05853140590135131
$ws.Range("G10").Value = 25.89747433333334
$ws.Range("H10").Value = 77.69242300000001
$ws.Range("I10").Value = 0.2930748622675039
$ws.Range("J10").Value = 0.2930748622675038
$ws.Range("M10").Value = 0.009290999999999999
$ws.Range("N10").Value = 0.027873
$ws.Range("O10").Value = 0.2547317242576836
$ws.Range("P10").Value = 0.2547317242576836
$ws.Range("Q10").Value = 0.240613434031
$ws.Range("R10").Value = 2.165520906279
$ws.Range("S10").Value = 0.0746554650019844
$ws.Range("T10").Value = 0.0746554650019844
$ws.Range("G11").Value = 20.774284
$ws.Range("H11").Value = 62.322852
$ws.Range("I11").Value = 0.2350970733145757
$ws.Range("J11").Value = 0.2350970733145757
$ws.Range("M11").Value = 0.01989833333333333
$ws.Range("N11").Value = 0.059695
$ws.Range("O11").Value = 0.5455534129646046
$ws.Range("P11").Value = 0.5455534129646046
$ws.Range("Q11").Value = 0.4133736277933333
$ws.Range("R11").Value = 3.72036265014
$ws.Range("S11").Value = 0.1282580107247566
$ws.Range("T11").Value = 0.1282580107247566
$ws.Range("G12").Value = 20.774284
$ws.Range("H12").Value = 62.322852
$ws.Range("I12").Value = 0.2350970733145757
$ws.Range("J12").Value = 0.2350970733145757
$ws.Range("O12").Value = 0.1997148627777118
$ws.Range("P12").Value = 0.1997148627777118
$ws.Range("Q12").Value = 0.1513268094173333
$ws.Range("R12").Value = 1.361941284756
$ws.Range("S12").Value = 0.04695237973646213
$ws.Range("T12").Value = 0.04695237973646213
$ws.Range("G13").Value = 20.774284
$ws.Range("H13").Value = 62.322852
$ws.Range("I13").Value = 0.2350970733145757
$ws.Range("J13").Value = 0.2350970733145757
$ws.Range("M13").Value = 0.009290999999999999
$ws.Range("N13").Value = 0.027873
$ws.Range("O13").Value = 0.2547317242576836
$ws.Range("P13").Value = 0.2547317242576836
$ws.Range("Q13").Value = 0.193013872644
$ws.Range("R13").Value = 1.737124853796
$ws.Range("S13").Value = 0.05988668285335692
$ws.Range("T13").Value = 0.05988668285335692
